# DialogueDB.xlsx update: add tutorial dialogue rows (events 튜토리얼1/2/3)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DialogueEntity")

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "튜토리얼1"
$ws.Range("C29").Value = "허수아비"
$ws.Range("D29").Value = "반가워! 오늘도 연습하러 왔네?`n다시 처음부터 연습해볼까?"
$ws.Range("D29").WrapText = $true
$ws.Range("E29").Value = $false

$ws.Range("A30").Value = 2
$ws.Range("B30").Value = "튜토리얼1"
$ws.Range("C30").Value = "허수아비"
$ws.Range("D30").Value = "먼저 움직이는 방법이야.`nWASD 키로 움직일 수 있어!"
$ws.Range("D30").WrapText = $true
$ws.Range("E30").Value = $false

$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "튜토리얼1"
$ws.Range("C31").Value = "허수아비"
$ws.Range("D31").Value = "공격은 마우스 왼클릭, 구르기는 마우스 우클릭이야.`n나를 한 번 공격해볼래?"
$ws.Range("D31").WrapText = $true
$ws.Range("E31").Value = $false

$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "튜토리얼1"
$ws.Range("C32").Value = "허수아비"
$ws.Range("D32").Value = "NPC나 오브젝트와 상호작용을 하고 싶다면 F키를 눌러봐"
$ws.Range("E32").Value = $false

$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "튜토리얼1"
$ws.Range("C33").Value = "허수아비"
$ws.Range("D33").Value = "잘했어!! 기본적인 조작은 모두 마스터했어.`n이제 화투에 대해 알려줄게!"
$ws.Range("D33").WrapText = $true
$ws.Range("E33").Value = $false

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "튜토리얼2"
$ws.Range("C34").Value = "허수아비"
$ws.Range("D34").Value = "화투는 용들이 아주 오래전부터 즐긴 전통놀이이자, 가장 좋아하는 놀이 중 하나야!!"
$ws.Range("E34").Value = $false

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "튜토리얼2"
$ws.Range("C35").Value = "허수아비"
$ws.Range("D35").Value = "용들은 화투 놀이를 즐기기도 했지만,`n마법의 모포를 이용해 화투에 힘을 담아 전투에 사용하기도 했어!"
$ws.Range("D35").WrapText = $true
$ws.Range("E35").Value = $false

$ws.Range("A36").Value = 8
$ws.Range("B36").Value = "튜토리얼2"
$ws.Range("C36").Value = "허수아비"
$ws.Range("D36").Value = "화투의 힘을 이용한 능력은 굉장했지..`n동물왕국을 모두 저지했을 정도였으니까!"
$ws.Range("D36").WrapText = $true
$ws.Range("E36").Value = $false

$ws.Range("A37").Value = 9
$ws.Range("B37").Value = "튜토리얼2"
$ws.Range("C37").Value = "허수아비"
$ws.Range("D37").Value = "하지만 이제는 화투 능력을 쓰는 용들이 거의 없는데…`n너는 아직도 쓰는구나??"
$ws.Range("D37").WrapText = $true
$ws.Range("E37").Value = $false

$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "튜토리얼2"
$ws.Range("C38").Value = "허수아비"
$ws.Range("D38").Value = "내가 화투 능력에 대해 자세히 알려주도록 할게!"
$ws.Range("E38").Value = $false

$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "튜토리얼3"
$ws.Range("C39").Value = "허수아비"
$ws.Range("D39").Value = "(뒤적뒤적..) 여기있다!`n이건 바로 화투 놀이중 하나인 섰다 족보책이야!`n이걸 받도록 해"
$ws.Range("D39").WrapText = $true
$ws.Range("E39").Value = $false

$ws.Range("A40").Value = 12
$ws.Range("B40").Value = "튜토리얼3"
$ws.Range("C40").Value = "허수아비"
$ws.Range("D40").Value = "섰다 족보책을 한 번 펼쳐볼래??`nK키를 누르면 책을 볼 수 있어."
$ws.Range("D40").WrapText = $true
$ws.Range("E40").Value = $false

$ws.Range("A41").Value = 13
$ws.Range("B41").Value = "튜토리얼3"
$ws.Range("C41").Value = "허수아비"
$ws.Range("D41").Value = "족보는 크게 광땡, 땡, 중간, 끗 순서대로야.`n각 카테고리를 살펴보면 어떤 화투패 조합이 있는지 알 수 있어."
$ws.Range("D41").WrapText = $true
$ws.Range("E41").Value = $false

$ws.Range("A42").Value = 14
$ws.Range("B42").Value = "튜토리얼3"
$ws.Range("C42").Value = "허수아비"
$ws.Range("D42").Value = "지금은 족보의 순서만 나와있지만 마법의 모포와 상호작용해서`n스킬을 얻으면 스킬 효과와 족보 시너지 효과를 알 수 있어!`n다양한 스킬을 얻고 사 용해보길 추천할께"
$ws.Range("D42").WrapText = $true
$ws.Range("E42").Value = $false

$ws.Range("A43").Value = 15
$ws.Range("B43").Value = "튜토리얼3"
$ws.Range("C43").Value = "허수아비"
$ws.Range("D43").Value = "임시로 한 가지 스킬만 알려주도록 할까? 이건 화투패 능력 중 하나야."
$ws.Range("E43").Value = $false

$ws.Range("A44").Value = 16
$ws.Range("B44").Value = "튜토리얼3"
$ws.Range("C44").Value = "허수아비"
$ws.Range("D44").Value = "Q키를 누르면 스킬을 사용할 수 있으니 한 번 사용해볼래?"
$ws.Range("E44").Value = $false

$ws.Range("A45").Value = 17
$ws.Range("B45").Value = "튜토리얼3"
$ws.Range("C45").Value = "허수아비"
$ws.Range("D45").Value = "잘했어!!`n이제 모든 훈련이 끝났으니 훈련장에서 나가봐도 좋아.`n이번엔 용과 시험에 꼭 통과하길 바라"
$ws.Range("D45").WrapText = $true
$ws.Range("E45").Value = $false

# Column D was manually resized (no longer best-fit/autosized)
$ws.Columns.Item(4).ColumnWidth = 66.77734375

# Final selection/zoom state as saved by the author
[void]$ws.Range("D41").Select()
$excel.ActiveWindow.Zoom = 78
